$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: change "Densidade" (Discrete) to "Furo" (Continuous)
$ws.Range("A2").Value = "Furo"
$ws.Range("B2").Value = 50
$ws.Range("C2").Value = 5
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = 0.95
$ws.Range("H2").Value = "Continuous"

# Row 3 (new): "Pino" (Continuous)
$ws.Range("A3").Value = "Pino"
$ws.Range("B3").Value = 40
$ws.Range("C3").Value = 3
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
# Trust level on row 3 is stored as text "0.95" (not a number) in target,
# so force text format before assigning to avoid auto-conversion to a number.
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "0.95"
$ws.Range("H3").Value = "Continuous"
